$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 765 (the existing
# "2026/12/29" block), shifting rows 765-806 down to 767-808.
$ws.Range("A765:A766").EntireRow.Insert()

# Force column A to be treated as text before writing, so the date-like
# string isn't auto-converted to a date serial number (matches the other
# date cells in this column, which are stored as plain text).
$ws.Range("A765:A766").NumberFormat = "@"

# Populate the two newly inserted rows with the new 2026/02/06 entries
# (continuing the same date block that ends at row 764).
$ws.Range("A765").Value = "2026/02/06"
$ws.Range("B765").Value = "金"
$ws.Range("C765").Value = 18
$ws.Range("D765").Value = 75

$ws.Range("A766").Value = "2026/02/06"
$ws.Range("B766").Value = "金"
$ws.Range("C766").Value = 22
$ws.Range("D766").Value = 78

# Drop the temporary text-number-format so the new rows end up with the
# same (unstyled) look as every other data row in the sheet.
$ws.Range("A765:D766").Style = "Normal"
